$wb = $excel.ActiveWorkbook

# --- survey sheet: add display.title.text / display.title.text.es columns ---
$survey = $wb.Worksheets.Item(1)

# Insert two new blank columns at F (pushes the existing prompt columns to H:I)
$survey.Columns.Item(6).Insert()
$survey.Columns.Item(6).Insert()

# Match the new title columns' width to the (now shifted) prompt.text column
$survey.Columns.Item(6).ColumnWidth = 23.833333333333332
$survey.Columns.Item(7).ColumnWidth = 23.833333333333332

# Header row
$survey.Range("F1").Value = "display.title.text"
$survey.Range("G1").Value = "display.title.text.es"
$survey.Range("I1").Value = "display.prompt.es"

# refrigerator_id row: title text
$survey.Range("F3").Value = "Refrigerator ID"
$survey.Range("G3").Value = "ID de Frigorífico"

# date_serviced row: title text
$survey.Range("F9").Value = "Date Serviced"
$survey.Range("G9").Value = "Fecha de Servicio"

# notes row: title text
$survey.Range("F10").Value = "Notes"
$survey.Range("G10").Value = "Notas"

# drop the trailing blank padding rows (12-21) that the sheet used to carry
$survey.Range("A12:A21").EntireRow.Delete()

# --- settings sheet: rename "spanish" labels to "es", bump form_version ---
$settings = $wb.Worksheets.Item(4)
$settings.Range("D1").Value = "display.title.text.es"
$settings.Range("F1").Value = "display.locale.text.es"
$settings.Range("A8").Value = "es"
$settings.Range("B5").Value = 20170717

# --- view state: make survey the active sheet/tab again ---
$settings.Activate()
$settings.Range("A9").Select()

$survey.Activate()
$survey.Range("G11").Select()
